# Atualizacao de bases das ligas, do dia: 20-06-2024 as 20:11
#
# The underlying match rows were re-sorted; each pair below swaps its
# entire data payload (columns B:AD) while leaving column A (the running
# rank number) untouched, matching the target workbook exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 83 and 84 (match ids 8075296 / 8075530)
$r1 = $ws.Range("B83:AD83")
$r2 = $ws.Range("B84:AD84")
$v1 = $r1.Value2
$v2 = $r2.Value2
$r1.Value2 = $v2
$r2.Value2 = $v1

# Swap rows 88 and 90 (match ids 8076438 / 8077795)
$r1 = $ws.Range("B88:AD88")
$r2 = $ws.Range("B90:AD90")
$v1 = $r1.Value2
$v2 = $r2.Value2
$r1.Value2 = $v2
$r2.Value2 = $v1

# Swap rows 97 and 98 (match ids 8158915 / 8163123)
$r1 = $ws.Range("B97:AD97")
$r2 = $ws.Range("B98:AD98")
$v1 = $r1.Value2
$v2 = $r2.Value2
$r1.Value2 = $v2
$r2.Value2 = $v1

# Swap rows 129 and 130 (match ids 8271342 / 8271343)
$r1 = $ws.Range("B129:AD129")
$r2 = $ws.Range("B130:AD130")
$v1 = $r1.Value2
$v2 = $r2.Value2
$r1.Value2 = $v2
$r2.Value2 = $v1
